# "Generate Report for Handoff"
#
# The localization-status report is being regenerated. For the four rows that
# were previously still queued at "low" priority (3b080f16, c91df214,
# e34b2465, ef37516b), the handoff step has now actually run:
#   - Priority flips from "low" to "ht" (matches the other already-handed-off rows)
#   - The zh-cn handoff timestamp advances from 16:36:17 to 16:36:40
#   - The shared "Latest HO Xliff Generate Date" / de-de handoff timestamp
#     advances from 16:36:21 to 16:36:46

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = 4, 5, 6, 7

foreach ($r in $rows) {
    # Overview sheet: Latest HO Xliff Generate Date (column G)
    $overview.Range("G$r").Value = "2016-08-19 16:36:46"

    # zh-cn sheet: Priority (column E) and Latest Handoff Datetime (column H)
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-19 16:36:40"

    # de-de sheet: Priority (column E) and Latest Handoff Datetime (column H)
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-19 16:36:46"
}
